$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update existing data values in column B (several entries -> 500)
# ---------------------------------------------------------------------
$ws.Range("B1").Value = 500
$ws.Range("B2").Value = 500
$ws.Range("B3").Value = 500
$ws.Range("B4").Value = 500
$ws.Range("B5").Value = 500
$ws.Range("B6").Value = 500
$ws.Range("B10").Value = 500
$ws.Range("B11").Value = 500
$ws.Range("B12").Value = 500
$ws.Range("B13").Value = 500

# ---------------------------------------------------------------------
# 2) Add four new rows (15-18): "No.15".."No.18" / 500
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "No.15"
$ws.Range("B15").Value = 500
$ws.Range("A16").Value = "No.16"
$ws.Range("B16").Value = 500
$ws.Range("A17").Value = "No.17"
$ws.Range("B17").Value = 500
$ws.Range("A18").Value = "No.18"
$ws.Range("B18").Value = 500

# ---------------------------------------------------------------------
# 3) Update the chart: title, series source range, overlap, position
# ---------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart

# Add chart title "Team Balance"
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Team Balance"

# Extend the series source ranges from row 14 to row 18
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(,Sheet1!`$A`$1:`$A`$18,Sheet1!`$B`$1:`$B`$18,1)"

# Remove the explicit bar overlap override
$chartGroup = $chart.ChartGroups(1)
$chartGroup.Overlap = $null

# Reposition / resize the chart on the worksheet
$co.Left = 261.95
$co.Top = 16.2
$co.Width = 790.525
$co.Height = 272.4

# ---------------------------------------------------------------------
# 4) Restore the cursor/selection like the saved workbook shows
# ---------------------------------------------------------------------
$ws.Range("D15").Select() | Out-Null
